$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update IsMember flags for existing rows 3 and 4 to TRUE
$ws.Range("D3").Value = $true
$ws.Range("D4").Value = $true

# Add new row 5: Id=4, Name="holy", PhoneNumber="123", IsMember=FALSE, CoffeeCount=0
$ws.Range("C5").NumberFormat = "@"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "holy"
$ws.Range("C5").Value = "123"
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = 0
$ws.Range("C5").Style = "Normal"

# Add new row 6: Id=5, Name="checkign price", PhoneNumber="69", IsMember=FALSE, CoffeeCount=0
$ws.Range("C6").NumberFormat = "@"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "checkign price"
$ws.Range("C6").Value = "69"
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = 0
$ws.Range("C6").Style = "Normal"

# Update selected cell to match the new active selection
$ws.Range("E4").Select()
